# "Add labels to the graph"
# Updates the GB (graph label) column on the Translation sheet with new
# axis-unit / tick labels, and appends four new rows of label text IDs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# Helper: write a value that must stay a genuine text cell even though it
# looks numeric (Excel/COM would otherwise auto-coerce "3.2" etc. into a
# numeric cell). Temporarily force text format, assign, then restore the
# cell to the default "Normal" style so no stray formatting is left behind.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Existing rows: refresh the graph's unit / value labels (column F).
Set-TextValue $ws.Range("F4") "[kHz]"
Set-TextValue $ws.Range("F5") "[ms]"
Set-TextValue $ws.Range("F8") "3.2"
Set-TextValue $ws.Range("F9") "0.8"
Set-TextValue $ws.Range("F10") "5"
Set-TextValue $ws.Range("F11") "10"

# New rows 12-15: additional single-use label text entries for the graph.
$ws.Range("B12").Value = "SingleUseId9"
$ws.Range("C12").Value = "Default"
$ws.Range("D12").Value = "Left"
$ws.Range("E12").Value = "LTR"
Set-TextValue $ws.Range("F12") "1.6"

$ws.Range("B13").Value = "SingleUseId10"
$ws.Range("C13").Value = "Default"
$ws.Range("D13").Value = "Left"
$ws.Range("E13").Value = "LTR"
Set-TextValue $ws.Range("F13") "2.4"

$ws.Range("B14").Value = "SingleUseId11"
$ws.Range("C14").Value = "Default"
$ws.Range("D14").Value = "Left"
$ws.Range("E14").Value = "LTR"
Set-TextValue $ws.Range("F14") "2.5"

$ws.Range("B15").Value = "SingleUseId12"
$ws.Range("C15").Value = "Default"
$ws.Range("D15").Value = "Left"
$ws.Range("E15").Value = "LTR"
Set-TextValue $ws.Range("F15") "7.5"
